$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1253.6923
$ws.Range("I70").Value = 1099.75
$ws.Range("K70").Value = 3299.25
$ws.Range("M70").Value = -3029.25
$ws.Range("H73").Value = 1253.6923
$ws.Range("I73").Value = 1099.75
$ws.Range("K73").Value = 3299.25
$ws.Range("M73").Value = -2363.25
$ws.Range("H97").Value = 19990
$ws.Range("J97").Value = 19990
$ws.Range("L97").Value = 59970
$ws.Range("N97").Value = -60962
$ws.Range("H112").Value = 1356.1086
$ws.Range("J112").Value = 1373.4318
$ws.Range("L112").Value = 4120.2954
$ws.Range("N112").Value = -6336.2954
$ws.Range("H137").Value = 3519.05
$ws.Range("I137").Value = 1338
$ws.Range("J137").Value = 3955.26
$ws.Range("K137").Value = 4014
$ws.Range("L137").Value = 11865.78
$ws.Range("M137").Value = -1464
$ws.Range("N137").Value = -16965.78
$ws.Range("H138").Value = 2058
$ws.Range("I138").Value = 1892.4286
$ws.Range("J138").Value = 2121.2183
$ws.Range("K138").Value = 5677.2858
$ws.Range("L138").Value = 6363.6549
$ws.Range("M138").Value = -537.2857999999997
$ws.Range("N138").Value = -16643.6549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28591.396
$ws.Range("I32").Value = 28249.05
$ws.Range("K32").Value = 28249.05
$ws.Range("M32").Value = -27962.05
$ws.Range("H61").Value = 2857.2092
$ws.Range("I61").Value = 1625.44
$ws.Range("J61").Value = 4568
$ws.Range("K61").Value = 1625.44
$ws.Range("L61").Value = 4568
$ws.Range("M61").Value = -1413.44
$ws.Range("N61").Value = -4992
$ws.Range("H74").Value = 1636.9231
$ws.Range("I74").Value = 1378.8572
$ws.Range("K74").Value = 1378.8572
$ws.Range("M74").Value = -504.8571999999999
$ws.Range("H77").Value = 1636.9231
$ws.Range("I77").Value = 1378.8572
$ws.Range("K77").Value = 6894.286
$ws.Range("M77").Value = -2526.286
$ws.Range("H122").Value = 3370.1538
$ws.Range("I122").Value = 3834.6667
$ws.Range("J122").Value = 2325
$ws.Range("K122").Value = 11504.0001
$ws.Range("L122").Value = 6975
$ws.Range("M122").Value = -9054.000100000001
$ws.Range("N122").Value = -11875
$ws.Range("H133").Value = 43298.875
$ws.Range("J133").Value = 43298.875
$ws.Range("L133").Value = 43298.875
$ws.Range("N133").Value = -48358.875
$ws.Range("H136").Value = 2857.2092
$ws.Range("I136").Value = 1625.44
$ws.Range("J136").Value = 4568
$ws.Range("K136").Value = 4876.32
$ws.Range("L136").Value = 13704
$ws.Range("M136").Value = -2326.32
$ws.Range("N136").Value = -18804

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1849.3914
$ws.Range("I86").Value = 1779.7778
$ws.Range("J86").Value = 2100
$ws.Range("K86").Value = 1779.7778
$ws.Range("L86").Value = 2100
$ws.Range("M86").Value = -656.7778000000001
$ws.Range("N86").Value = -4346
$ws.Range("H89").Value = 1849.3914
$ws.Range("I89").Value = 1779.7778
$ws.Range("J89").Value = 2100
$ws.Range("K89").Value = 8898.889000000001
$ws.Range("L89").Value = 10500
$ws.Range("M89").Value = -3282.889000000001
$ws.Range("N89").Value = -21732
$ws.Range("H137").Value = 40373
$ws.Range("J137").Value = 40373
$ws.Range("L137").Value = 40373
$ws.Range("N137").Value = -50573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 41953.4
$ws.Range("J59").Value = 41953.4
$ws.Range("L59").Value = 41953.4
$ws.Range("N59").Value = -44243.4
$ws.Range("H62").Value = 2998.5454
$ws.Range("I62").Value = 2998.8
$ws.Range("J62").Value = 2998.3333
$ws.Range("K62").Value = 2998.8
$ws.Range("L62").Value = 2998.3333
$ws.Range("M62").Value = -2374.8
$ws.Range("N62").Value = -4246.3333
$ws.Range("H65").Value = 2998.5454
$ws.Range("I65").Value = 2998.8
$ws.Range("J65").Value = 2998.3333
$ws.Range("K65").Value = 14994
$ws.Range("L65").Value = 14991.6665
$ws.Range("M65").Value = -11874
$ws.Range("N65").Value = -21231.6665
$ws.Range("H68").Value = 56999.5
$ws.Range("J68").Value = 56999.5
$ws.Range("L68").Value = 56999.5
$ws.Range("N68").Value = -58497.5
$ws.Range("H71").Value = 56999.5
$ws.Range("J71").Value = 56999.5
$ws.Range("L71").Value = 170998.5
$ws.Range("N71").Value = -178486.5
$ws.Range("H74").Value = 14000
$ws.Range("J74").Value = 14000
$ws.Range("L74").Value = 14000
$ws.Range("N74").Value = -15748
$ws.Range("H77").Value = 14000
$ws.Range("J77").Value = 14000
$ws.Range("L77").Value = 42000
$ws.Range("N77").Value = -50736
$ws.Range("H132").Value = 42170.086
$ws.Range("I132").Value = 1617
$ws.Range("J132").Value = 179036.75
$ws.Range("K132").Value = 4851
$ws.Range("L132").Value = 537110.25
$ws.Range("M132").Value = -2321
$ws.Range("N132").Value = -542170.25
$ws.Range("H134").Value = 609872.4399999999
$ws.Range("I134").Value = 880.5294
$ws.Range("K134").Value = 2641.5882
$ws.Range("M134").Value = -106.5882000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 64877.453
$ws.Range("J131").Value = 88477.87
$ws.Range("L131").Value = 265433.61
$ws.Range("N131").Value = -275513.61

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1911.1111
$ws.Range("I122").Value = 1880
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5640
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3190
$ws.Range("N122").Value = -10750
$ws.Range("H123").Value = 12654.5
$ws.Range("J123").Value = 12654.5
$ws.Range("L123").Value = 12654.5
$ws.Range("N123").Value = -17554.5
$ws.Range("H126").Value = 2384.111
$ws.Range("I126").Value = 3145.6667
$ws.Range("J126").Value = 2003.3334
$ws.Range("K126").Value = 9437.000100000001
$ws.Range("L126").Value = 6010.0002
$ws.Range("M126").Value = -6967.000100000001
$ws.Range("N126").Value = -10950.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3508.7144
$ws.Range("I9").Value = 387.5
$ws.Range("J9").Value = 7670.3335
$ws.Range("K9").Value = 387.5
$ws.Range("L9").Value = 7670.3335
$ws.Range("M9").Value = -163.5
$ws.Range("N9").Value = -8118.3335
$ws.Range("H68").Value = 2833.2222
$ws.Range("I68").Value = 2099.75
$ws.Range("J68").Value = 3420
$ws.Range("K68").Value = 2099.75
$ws.Range("L68").Value = 3420
$ws.Range("M68").Value = -1350.75
$ws.Range("N68").Value = -4918
$ws.Range("H71").Value = 2833.2222
$ws.Range("I71").Value = 2099.75
$ws.Range("J71").Value = 3420
$ws.Range("K71").Value = 10498.75
$ws.Range("L71").Value = 17100
$ws.Range("M71").Value = -6754.75
$ws.Range("N71").Value = -24588

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 83421
$ws.Range("J133").Value = 83421
$ws.Range("L133").Value = 83421
$ws.Range("N133").Value = -93541
$ws.Range("H136").Value = 22684.834
$ws.Range("I136").Value = 77954
$ws.Range("J136").Value = 2156.2856
$ws.Range("K136").Value = 233862
$ws.Range("L136").Value = 6468.8568
$ws.Range("M136").Value = -231312
$ws.Range("N136").Value = -11568.8568
